$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the value previously in E42 (last row before the new quarter is added)
$ws.Range("E42").Value = 7215

# Append the new quarterly row (43) with the new period "01-04-2021".
# The leading apostrophe forces Excel to store the value as literal text
# instead of auto-converting the dd-mm-yyyy-looking string into a date
# serial number (matching how the other period labels in column A are
# stored as plain shared strings).
$ws.Range("A43").Value = "'01-04-2021"
$ws.Range("A43").Style = "Normal"

$ws.Range("B43").Value = 4580
$ws.Range("C43").Value = 2332
$ws.Range("D43").Value = 1981
$ws.Range("E43").Value = 7418
